$d = $word.ActiveDocument

$find = "2022 Ημερομηνίες παρατήρησης για τον αστερισμό του Αστερισμός Λέων: 14-23 Απριλίου, 14-23 Μαΐου"
$replace = "2022 Ημερομηνίες παρατήρησης για τον  Αστερισμός Λέων: 14-23 Απριλίου, 14-23 Μαΐου"

$range = $d.Content
$range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
